$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 264

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

for ($r = 2; $r -le 11; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2

    $ws.Range("S$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALINGSAS/artfynd/' + $label + '.xlsx", "' + $label + '")'
    $ws.Range("T$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALINGSAS/kartor/' + $label + '.png", "' + $label + '")'
    $ws.Range("V$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALINGSAS/klagomål/' + $label + '.docx", "' + $label + '")'
    $ws.Range("W$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALINGSAS/klagomålsmail/' + $label + '.docx", "' + $label + '")'
    $ws.Range("X$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALINGSAS/tillsyn/' + $label + '.docx", "' + $label + '")'
    $ws.Range("Y$r").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALINGSAS/tillsynsmail/' + $label + '.docx", "' + $label + '")'
}
